$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume updates as scraped on Tue Jul  2 11:50:38 UTC 2024

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.749.07'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.450.97'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.70'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.52'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '8.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.415'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.042.96'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.29'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.442.07'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.786.81'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.40'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.61'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.01'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '387.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.568'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '75.32'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.585.76'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('E27').Value = '  +2.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.70'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.01'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.24'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('E35').Value = '  +2.68%  '
$ws.Range('E36').Value = '  +4.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '31.87'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.94'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '169.17'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.485.53'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.84'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.567.75'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.91'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.24'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('E50').Value = '  -3.06%  '
$ws.Range('E51').Value = '  -0.09%  '
